$p = $ppt.ActivePresentation

# Delete slide 3 first so slide 2's index isn't shifted while removing it.
$p.Slides.Item(3).Delete()
$p.Slides.Item(2).Delete()
